$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: fill in the "punch.mp3" sound-effect entry ---
$ws.Range("A33").Value = "punch.mp3"
$ws.Range("B33").Value = "Sound of a punch"
$ws.Range("C33").Value = "1 sec"

# D33 gets a hyperlink (reuse the existing Hyperlink cell style/border by
# copying the format from another hyperlink cell, then set the real value
# and hyperlink separately so the style/text stay correct).
$ws.Range("D9").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = "https://www.youtube.com/watch?v=RHFN4-BLcIo"
$ws.Hyperlinks.Add($ws.Range("D33"), "https://www.youtube.com/watch?v=RHFN4-BLcIo")

$ws.Range("E33").Value = "Done"

# --- Row 34 (new row): "memento.mp3" main menu music entry ---
$ws.Range("A34").Value = "memento.mp3"
$ws.Range("B34").Value = "Main Menu Music"
$ws.Range("C34").Value = "3min:31sec"
$ws.Range("D34").Value = "Composer: Myuu"
$ws.Range("E34").Value = "Done"

# --- Row 19: "main_menu.png" picture entry ---
$ws.Range("A19").Value = "main_menu.png"
$ws.Range("B19").Value = "Main Menu Picture"
$ws.Range("C19").Value = "1280x1080px"
$ws.Range("E19").Value = "Placeholder"

# Update the selection shown when the workbook is reopened.
$ws.Range("E19").Select()
